$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 6217
$ws.Range("I131").Value = 3028.3333
$ws.Range("J131").Value = 11000
$ws.Range("K131").Value = 9084.999899999999
$ws.Range("L131").Value = 33000
$ws.Range("M131").Value = -4044.999899999999
$ws.Range("N131").Value = -43080

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2974.25
$ws.Range("I2").Value = 1313
$ws.Range("J2").Value = 5300
$ws.Range("K2").Value = 1313
$ws.Range("L2").Value = 5300
$ws.Range("M2").Value = -1200
$ws.Range("N2").Value = -5526

$ws.Range("H74").Value = 1487.125
$ws.Range("I74").Value = 1035.8182
$ws.Range("J74").Value = 2480
$ws.Range("K74").Value = 1035.8182
$ws.Range("L74").Value = 2480
$ws.Range("M74").Value = -161.8181999999999
$ws.Range("N74").Value = -4228

$ws.Range("H77").Value = 1487.125
$ws.Range("I77").Value = 1035.8182
$ws.Range("J77").Value = 2480
$ws.Range("K77").Value = 5179.090999999999
$ws.Range("L77").Value = 12400
$ws.Range("M77").Value = -811.0909999999994
$ws.Range("N77").Value = -21136

$ws.Range("H116").Value = 2974.25
$ws.Range("I116").Value = 1313
$ws.Range("J116").Value = 5300
$ws.Range("K116").Value = 1313
$ws.Range("L116").Value = 5300
$ws.Range("M116").Value = 981
$ws.Range("N116").Value = -9888

$ws.Range("H122").Value = 2253.5854
$ws.Range("I122").Value = 2136.7896
$ws.Range("J122").Value = 3733
$ws.Range("K122").Value = 6410.3688
$ws.Range("L122").Value = 11199
$ws.Range("M122").Value = -3960.3688
$ws.Range("N122").Value = -16099

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2974.25
$ws.Range("I3").Value = 1313
$ws.Range("J3").Value = 5300
$ws.Range("K3").Value = 1313
$ws.Range("L3").Value = 5300
$ws.Range("M3").Value = -1199
$ws.Range("N3").Value = -5528

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 1171.4286
$ws.Range("I62").Value = 1033.3334
$ws.Range("J62").Value = 2000
$ws.Range("K62").Value = 1033.3334
$ws.Range("L62").Value = 2000
$ws.Range("M62").Value = -409.3334
$ws.Range("N62").Value = -3248

$ws.Range("H65").Value = 1171.4286
$ws.Range("I65").Value = 1033.3334
$ws.Range("J65").Value = 2000
$ws.Range("K65").Value = 5166.666999999999
$ws.Range("L65").Value = 10000
$ws.Range("M65").Value = -2046.666999999999
$ws.Range("N65").Value = -16240

$ws.Range("H99").Value = 2127.1667
$ws.Range("I99").Value = 2087.0334
$ws.Range("J99").Value = 2327.8333
$ws.Range("K99").Value = 2087.0334
$ws.Range("L99").Value = 2327.8333
$ws.Range("M99").Value = -589.0333999999998
$ws.Range("N99").Value = -5323.8333

$ws.Range("H126").Value = 2127.1667
$ws.Range("I126").Value = 2087.0334
$ws.Range("J126").Value = 2327.8333
$ws.Range("K126").Value = 6261.100199999999
$ws.Range("L126").Value = 6983.499899999999
$ws.Range("M126").Value = -3791.100199999999
$ws.Range("N126").Value = -11923.4999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 322.14285
$ws.Range("I2").Value = 40
$ws.Range("J2").Value = 343.84616
$ws.Range("K2").Value = 240
$ws.Range("L2").Value = 2063.07696
$ws.Range("M2").Value = -127
$ws.Range("N2").Value = -2289.07696

$ws.Range("H22").Value = 55556860
$ws.Range("I22").Value = 1348
$ws.Range("J22").Value = 62501296
$ws.Range("K22").Value = 4044
$ws.Range("L22").Value = 187503888
$ws.Range("M22").Value = -3875
$ws.Range("N22").Value = -187504226

$ws.Range("H27").Value = 55556860
$ws.Range("I27").Value = 1348
$ws.Range("J27").Value = 62501296
$ws.Range("K27").Value = 4044
$ws.Range("L27").Value = 187503888
$ws.Range("M27").Value = -3942
$ws.Range("N27").Value = -187504092

$ws.Range("H33").Value = 241.19048
$ws.Range("I33").Value = 57
$ws.Range("J33").Value = 540.5
$ws.Range("K33").Value = 342
$ws.Range("L33").Value = 3243
$ws.Range("M33").Value = -59
$ws.Range("N33").Value = -3809

$ws.Range("H39").Value = 1653.125
$ws.Range("I39").Value = 500
$ws.Range("J39").Value = 3135.7144
$ws.Range("K39").Value = 1500
$ws.Range("L39").Value = 9407.143199999999
$ws.Range("M39").Value = -1206
$ws.Range("N39").Value = -9995.143199999999

$ws.Range("H44").Value = 897.9
$ws.Range("I44").Value = 499.5
$ws.Range("J44").Value = 997.5
$ws.Range("K44").Value = 1498.5
$ws.Range("L44").Value = 2992.5
$ws.Range("M44").Value = -1100.5
$ws.Range("N44").Value = -3788.5

$ws.Range("H49").Value = 1873.3334
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 1873.3334
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 5620.0002
$ws.Range("N49").Value = -5932.0002

$ws.Range("H58").Value = 2743.6365
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 2743.6365
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 8230.9095
$ws.Range("N58").Value = -8486.9095

$ws.Range("H70").Value = 1354.4445
$ws.Range("I70").Value = 1031.6666
$ws.Range("J70").Value = 2000
$ws.Range("K70").Value = 3094.9998
$ws.Range("L70").Value = 6000
$ws.Range("M70").Value = -2779.9998
$ws.Range("N70").Value = -6630

$ws.Range("H73").Value = 1354.4445
$ws.Range("I73").Value = 1031.6666
$ws.Range("J73").Value = 2000
$ws.Range("K73").Value = 3094.9998
$ws.Range("L73").Value = 6000
$ws.Range("M73").Value = -2002.9998
$ws.Range("N73").Value = -8184

$ws.Range("H76").Value = 4809.524
$ws.Range("I76").Value = 1000
$ws.Range("J76").Value = 5000
$ws.Range("K76").Value = 3000
$ws.Range("L76").Value = 15000
$ws.Range("M76").Value = -2617
$ws.Range("N76").Value = -15766

$ws.Range("H79").Value = 4809.524
$ws.Range("I79").Value = 1000
$ws.Range("J79").Value = 5000
$ws.Range("K79").Value = 3000
$ws.Range("L79").Value = 15000
$ws.Range("M79").Value = -1674
$ws.Range("N79").Value = -17652

$ws.Range("H110").Value = 6031.6665
$ws.Range("I110").Value = 2033.3334
$ws.Range("J110").Value = 10030
$ws.Range("K110").Value = 6100.0002
$ws.Range("L110").Value = 30090
$ws.Range("M110").Value = -2010.0002
$ws.Range("N110").Value = -38270

$ws.Range("H111").Value = 708.5
$ws.Range("I111").Value = 708.5
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 2125.5
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = 941.5

$ws.Range("H112").Value = 393249.22
$ws.Range("I112").Value = 3356243.2
$ws.Range("J112").Value = 3381.5789
$ws.Range("K112").Value = 10068729.6
$ws.Range("L112").Value = 10144.7367
$ws.Range("M112").Value = -10067621.6
$ws.Range("N112").Value = -12360.7367

$ws.Range("H124").Value = 4983.3335
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 4983.3335
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 14950.0005
$ws.Range("N124").Value = -24770.0005
$ws.Range("M124").ClearContents()

$ws.Range("H125").Value = 4343.75
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 4343.75
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 13031.25
$ws.Range("N125").Value = -22871.25
$ws.Range("M125").ClearContents()

$ws.Range("H126").Value = 1745
$ws.Range("I126").Value = 1000
$ws.Range("J126").Value = 2490
$ws.Range("K126").Value = 3000
$ws.Range("L126").Value = 7470
$ws.Range("M126").Value = 1940
$ws.Range("N126").Value = -17350

$ws.Range("H129").Value = 27780020
$ws.Range("I129").Value = 4800
$ws.Range("J129").Value = 37038428
$ws.Range("K129").Value = 14400
$ws.Range("L129").Value = 111115284
$ws.Range("M129").Value = -9400
$ws.Range("N129").Value = -111125284

$ws.Range("H130").Value = 0
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("M130").ClearContents()

$ws.Range("H131").Value = 730.4792
$ws.Range("I131").Value = 423.07693
$ws.Range("J131").Value = 778.6265
$ws.Range("K131").Value = 1269.23079
$ws.Range("L131").Value = 2335.8795
$ws.Range("M131").Value = 3770.76921
$ws.Range("N131").Value = -12415.8795

$ws.Range("H133").Value = 7300.9
$ws.Range("I133").Value = 3810.5557
$ws.Range("J133").Value = 8067.073
$ws.Range("K133").Value = 11431.6671
$ws.Range("L133").Value = 24201.219
$ws.Range("M133").Value = -6371.667099999999
$ws.Range("N133").Value = -34321.219

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 35714956
$ws.Range("I113").Value = 495.83334
$ws.Range("J113").Value = 250001710
$ws.Range("K113").Value = 1487.50002
$ws.Range("L113").Value = 750005130
$ws.Range("M113").Value = 682.4999800000001
$ws.Range("N113").Value = -750009470
